# QB Website updated 12/13
# Append the newly-added occurrence numbers (83, 87) to the "Occurrence"
# column (column E) entries on Sheet1. The edit touches rows 2 and 4-9
# (row 3 / cell E3 is left unchanged), and moves the active selection
# from E9 to E4.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Order matters here only insofar as it determines the order new shared
# strings get appended to the workbook's string table; it mirrors the
# order the numbers were actually edited in upstream.
$ws.Range("E5").Value = "1 to 21, 23, 25, 27, 35, 39, 43, 47, 51, 55, 59, 67, 71, 75, 83, 87"
$ws.Range("E6").Value = "19, 21, 23, 25, 27, 35, 39, 43, 47, 51, 55,  59, 67, 71, 75, 83, 87"
$ws.Range("E2").Value = "19, 21, 23, 25, 27, 35, 39, 43, 47, 51, 55, 59, 67, 71, 75, 83, 87"
$ws.Range("E7").Value = "19, 21, 23, 25, 27, 35, 39, 43, 47, 51, 55, 59, 67, 71, 75, 83, 87"
$ws.Range("E8").Value = "21, 23, 25, 27, 35, 39, 43, 47, 51, 55, 59, 67, 71, 75, 83, 87"
$ws.Range("E9").Value = "21, 23, 25, 27, 35, 39, 43, 47, 51, 55, 59, 67, 71, 75, 83, 87"
$ws.Range("E4").Value = "51, 55, 59, 67, 71, 75, 83, 87"

# Move the selection/active cell to E4 (was E9), and scroll so the
# frozen/top-left cell resets to the sheet's natural top-left instead of A5.
[void]$ws.Range("E4").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
